$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "HMP Fred"
$ws.Range("C2").Value = "Freds County Court"

$ws.Range("C2").Font.Name = "Arial"
$ws.Range("C2").Font.Size = 10
$ws.Range("C2").Font.Color = 0

$ws.Rows.Item(2).AutoFit()

$ws.Range("B2").Select()
